# Apply weekly crypto price/volume refresh (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.899.27'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.903.92'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7998'
$ws.Range('E5').Value = '  +5.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.20'
$ws.Range('E6').Value = '  +0.39%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3122'
$ws.Range('E8').Value = '  +2.40%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.18'
$ws.Range('E9').Value = '  +3.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06871'
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07981'
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '1.913.01'
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7343'
$ws.Range('E13').Value = '  -2.22%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.166'
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.67'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '29.907.61'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.89'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.849'
$ws.Range('E18').Value = '  -1.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.86'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007690'
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').Value = '2.155.38'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.902'
$ws.Range('E24').Value = '  -1.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '168.07'
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.179'
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1412'
$ws.Range('E27').Value = '  +8.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.84'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.013'
$ws.Range('E29').Value = '  -0.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.359'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.514'
$ws.Range('E31').Value = '  -0.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.280'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05529'
$ws.Range('E33').Value = '  +3.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.056'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.254'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7282'
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.727'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01923'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.791'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.129'
$ws.Range('E40').Value = '  -1.11%  '
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '71.88'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.002'
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8354'
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.862'
$ws.Range('E45').Value = '  -2.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.55'
$ws.Range('E46').Value = '  -0.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.535'
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.707'
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.063.18'
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '976.98'
$ws.Range('E50').Value = '  +6.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.15'
$ws.Range('E51').Value = '  -0.05%  '
